$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing rows so the new row matches the
# established style pattern (A: date style, B: decimal style used in
# row 3, C: wrapped text style), without introducing new style entries.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)

# Fill in the new log entry values.
$ws.Range("A9").Value = 43894
$ws.Range("B9").Value = 2.25
$ws.Range("C9").Value = "Absence pour le recrutement"

# Move the active selection to the newly added cell, matching the
# workbook's saved cursor position.
$ws.Range("C9").Select()
